$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = [double]"0.07602548599243164"
$ws.Cells.Item(2, 2).Value = [double]"0.9771606922149658"
$ws.Cells.Item(2, 3).Value = [double]"0.07422774285078049"
$ws.Cells.Item(2, 4).Value = [double]"0.9905523061752319"
$ws.Cells.Item(3, 1).Value = [double]"0.01085434574633837"
$ws.Cells.Item(3, 2).Value = [double]"0.9985154271125793"
$ws.Cells.Item(3, 3).Value = [double]"0.06154094636440277"
$ws.Cells.Item(3, 4).Value = [double]"0.9905523061752319"
$ws.Cells.Item(4, 1).Value = [double]"0.007147168274968863"
$ws.Cells.Item(4, 2).Value = [double]"0.998629629611969"
$ws.Cells.Item(4, 3).Value = [double]"0.03475822880864143"
$ws.Cells.Item(4, 4).Value = [double]"0.9934592843055725"
$ws.Cells.Item(5, 1).Value = [double]"0.004885368049144745"
$ws.Cells.Item(5, 2).Value = [double]"0.9989722371101379"
$ws.Cells.Item(5, 3).Value = [double]"0.007020926102995872"
$ws.Cells.Item(5, 4).Value = [double]"0.9978197813034058"
$ws.Cells.Item(6, 1).Value = [double]"0.001769306254573166"
$ws.Cells.Item(6, 2).Value = [double]"0.9995431900024414"
$ws.Cells.Item(6, 3).Value = [double]"0.009036371484398842"
$ws.Cells.Item(6, 4).Value = [double]"0.9985465407371521"
$ws.Cells.Item(7, 1).Value = [double]"0.002138434210792184"
$ws.Cells.Item(7, 2).Value = [double]"0.9995431900024414"
$ws.Cells.Item(7, 3).Value = [double]"0.00184949766844511"
$ws.Cells.Item(7, 4).Value = [double]"1"
$ws.Cells.Item(8, 1).Value = [double]"0.001323148841038346"
$ws.Cells.Item(8, 2).Value = [double]"0.9997525811195374"
$ws.Cells.Item(8, 3).Value = [double]"0.00133598071988672"
$ws.Cells.Item(8, 4).Value = [double]"0.9992732405662537"
$ws.Cells.Item(9, 1).Value = [double]"0.001217141398228705"
$ws.Cells.Item(9, 2).Value = [double]"0.9996764659881592"
$ws.Cells.Item(9, 3).Value = [double]"7.343779725488275E-05"
$ws.Cells.Item(9, 4).Value = [double]"1"
$ws.Cells.Item(10, 1).Value = [double]"0.0007175981299951673"
$ws.Cells.Item(10, 2).Value = [double]"0.9998286962509155"
$ws.Cells.Item(10, 3).Value = [double]"0.000195938817341812"
$ws.Cells.Item(10, 4).Value = [double]"1"
$ws.Cells.Item(11, 1).Value = [double]"0.001031169202178717"
$ws.Cells.Item(11, 2).Value = [double]"0.9997715950012207"
$ws.Cells.Item(11, 3).Value = [double]"2.48798151005758E-05"
$ws.Cells.Item(11, 4).Value = [double]"1"
$ws.Cells.Item(12, 1).Value = [double]"0.0007697513210587204"
$ws.Cells.Item(12, 2).Value = [double]"0.9998096823692322"
$ws.Cells.Item(12, 3).Value = [double]"4.395343785290606E-05"
$ws.Cells.Item(12, 4).Value = [double]"1"
$ws.Cells.Item(13, 1).Value = [double]"0.0005132107180543244"
$ws.Cells.Item(13, 2).Value = [double]"0.9998857975006104"
$ws.Cells.Item(13, 3).Value = [double]"6.74099328534794E-06"
$ws.Cells.Item(13, 4).Value = [double]"1"
$ws.Cells.Item(14, 1).Value = [double]"0.0006111887050792575"
$ws.Cells.Item(14, 2).Value = [double]"0.9998286962509155"
$ws.Cells.Item(14, 3).Value = [double]"0.004253774415701628"
$ws.Cells.Item(14, 4).Value = [double]"0.9985465407371521"
$ws.Cells.Item(15, 1).Value = [double]"0.001032147789373994"
$ws.Cells.Item(15, 2).Value = [double]"0.9997144937515259"
$ws.Cells.Item(15, 3).Value = [double]"1.950105070136487E-05"
$ws.Cells.Item(15, 4).Value = [double]"1"
$ws.Cells.Item(16, 1).Value = [double]"0.0003594239242374897"
$ws.Cells.Item(16, 2).Value = [double]"0.9999048113822937"
$ws.Cells.Item(16, 3).Value = [double]"1.457721054975991E-06"
$ws.Cells.Item(16, 4).Value = [double]"1"
$ws.Cells.Item(17, 1).Value = [double]"0.0006362181156873703"
$ws.Cells.Item(17, 2).Value = [double]"0.9997715950012207"
$ws.Cells.Item(17, 3).Value = [double]"1.725155743770301E-05"
$ws.Cells.Item(17, 4).Value = [double]"1"
$ws.Cells.Item(18, 1).Value = [double]"0.0001598382368683815"
$ws.Cells.Item(18, 2).Value = [double]"0.9999428987503052"
$ws.Cells.Item(18, 3).Value = [double]"3.314510649943259E-06"
$ws.Cells.Item(18, 4).Value = [double]"1"
$ws.Cells.Item(19, 1).Value = [double]"0.0006894905818626285"
$ws.Cells.Item(19, 2).Value = [double]"0.9998857975006104"
$ws.Cells.Item(19, 3).Value = [double]"6.918030521774199E-06"
$ws.Cells.Item(19, 4).Value = [double]"1"
$ws.Cells.Item(20, 1).Value = [double]"0.00113664660602808"
$ws.Cells.Item(20, 2).Value = [double]"0.9998096823692322"
$ws.Cells.Item(20, 3).Value = [double]"3.852990630548447E-06"
$ws.Cells.Item(20, 4).Value = [double]"1"
$ws.Cells.Item(21, 1).Value = [double]"9.051190863829106E-05"
$ws.Cells.Item(21, 2).Value = [double]"0.9999619126319885"
$ws.Cells.Item(21, 3).Value = [double]"1.272373992833309E-06"
$ws.Cells.Item(21, 4).Value = [double]"1"
$ws.Cells.Item(22, 1).Value = [double]"0.0002641715691424906"
$ws.Cells.Item(22, 2).Value = [double]"0.9999809861183167"
$ws.Cells.Item(22, 3).Value = [double]"2.372712060605409E-06"
$ws.Cells.Item(22, 4).Value = [double]"1"
$ws.Cells.Item(23, 1).Value = [double]"0.0005599940195679665"
$ws.Cells.Item(23, 2).Value = [double]"0.999866783618927"
$ws.Cells.Item(23, 3).Value = [double]"4.178047220193548E-06"
$ws.Cells.Item(23, 4).Value = [double]"1"
$ws.Cells.Item(24, 1).Value = [double]"0.000627832196187228"
$ws.Cells.Item(24, 2).Value = [double]"0.9998857975006104"
$ws.Cells.Item(24, 3).Value = [double]"8.909722964745015E-06"
$ws.Cells.Item(24, 4).Value = [double]"1"
$ws.Cells.Item(25, 1).Value = [double]"0.0002438678930047899"
$ws.Cells.Item(25, 2).Value = [double]"0.9999428987503052"
$ws.Cells.Item(25, 3).Value = [double]"8.888677984941751E-05"
$ws.Cells.Item(25, 4).Value = [double]"1"
$ws.Cells.Item(26, 1).Value = [double]"0.0009342542034573853"
$ws.Cells.Item(26, 2).Value = [double]"0.9997715950012207"
$ws.Cells.Item(26, 3).Value = [double]"2.084785410261247E-05"
$ws.Cells.Item(26, 4).Value = [double]"1"
$ws.Cells.Item(27, 1).Value = [double]"0.0001685831375652924"
$ws.Cells.Item(27, 2).Value = [double]"0.9999238848686218"
$ws.Cells.Item(27, 3).Value = [double]"1.92800671356963E-05"
$ws.Cells.Item(27, 4).Value = [double]"1"
$ws.Cells.Item(28, 1).Value = [double]"0.0003457276325207204"
$ws.Cells.Item(28, 2).Value = [double]"0.9998857975006104"
$ws.Cells.Item(28, 3).Value = [double]"0.0001156580838141963"
$ws.Cells.Item(28, 4).Value = [double]"1"
$ws.Cells.Item(29, 1).Value = [double]"0.0001903246011352167"
$ws.Cells.Item(29, 2).Value = [double]"0.9999428987503052"
$ws.Cells.Item(29, 3).Value = [double]"5.260014222585596E-05"
$ws.Cells.Item(29, 4).Value = [double]"1"
$ws.Cells.Item(30, 1).Value = [double]"0.000285543326754123"
$ws.Cells.Item(30, 2).Value = [double]"0.9999048113822937"
$ws.Cells.Item(30, 3).Value = [double]"0.00164121063426137"
$ws.Cells.Item(30, 4).Value = [double]"0.9992732405662537"
$ws.Cells.Item(31, 1).Value = [double]"7.938235648907721E-05"
$ws.Cells.Item(31, 2).Value = [double]"0.9999809861183167"
$ws.Cells.Item(31, 3).Value = [double]"1.75736549863359E-05"
$ws.Cells.Item(31, 4).Value = [double]"1"
$ws.Cells.Item(32, 1).Value = [double]"0.0004082358500454575"
$ws.Cells.Item(32, 2).Value = [double]"0.9998857975006104"
$ws.Cells.Item(32, 3).Value = [double]"3.535087898853817E-06"
$ws.Cells.Item(32, 4).Value = [double]"1"
$ws.Cells.Item(33, 1).Value = [double]"0.0006414633826352656"
$ws.Cells.Item(33, 2).Value = [double]"0.9999048113822937"
$ws.Cells.Item(33, 3).Value = [double]"1.747312921906996E-06"
$ws.Cells.Item(33, 4).Value = [double]"1"
$ws.Cells.Item(34, 1).Value = [double]"0.0007026436505839229"
$ws.Cells.Item(34, 2).Value = [double]"0.9998857975006104"
$ws.Cells.Item(34, 3).Value = [double]"3.376153927092673E-06"
$ws.Cells.Item(34, 4).Value = [double]"1"
$ws.Cells.Item(35, 1).Value = [double]"0.0006924913031980395"
$ws.Cells.Item(35, 2).Value = [double]"0.9998477101325989"
$ws.Cells.Item(35, 3).Value = [double]"0.0001846192753873765"
$ws.Cells.Item(35, 4).Value = [double]"1"
$ws.Cells.Item(36, 1).Value = [double]"0.0002042878477368504"
$ws.Cells.Item(36, 2).Value = [double]"0.9999428987503052"
$ws.Cells.Item(36, 3).Value = [double]"0.0001047378609655425"
$ws.Cells.Item(36, 4).Value = [double]"1"
$ws.Cells.Item(37, 1).Value = [double]"7.083082891767845E-05"
$ws.Cells.Item(37, 2).Value = [double]"0.9999809861183167"
$ws.Cells.Item(37, 3).Value = [double]"4.839326265937416E-06"
$ws.Cells.Item(37, 4).Value = [double]"1"
$ws.Cells.Item(38, 1).Value = [double]"2.744052835623734E-05"
$ws.Cells.Item(38, 2).Value = [double]"1"
$ws.Cells.Item(38, 3).Value = [double]"7.846702487768198E-07"
$ws.Cells.Item(38, 4).Value = [double]"1"
$ws.Cells.Item(39, 1).Value = [double]"1.691205943643581E-05"
$ws.Cells.Item(39, 2).Value = [double]"1"
$ws.Cells.Item(39, 3).Value = [double]"1.206733628578149E-07"
$ws.Cells.Item(39, 4).Value = [double]"1"
$ws.Cells.Item(40, 1).Value = [double]"0.0002056228258879855"
$ws.Cells.Item(40, 2).Value = [double]"0.9999238848686218"
$ws.Cells.Item(40, 3).Value = [double]"0.0006103392806835473"
$ws.Cells.Item(40, 4).Value = [double]"0.9992732405662537"
$ws.Cells.Item(41, 1).Value = [double]"0.0004009988042525947"
$ws.Cells.Item(41, 2).Value = [double]"0.9999048113822937"
$ws.Cells.Item(41, 3).Value = [double]"3.483045475149993E-06"
$ws.Cells.Item(41, 4).Value = [double]"1"
$ws.Cells.Item(42, 1).Value = [double]"0.0002500255068298429"
$ws.Cells.Item(42, 2).Value = [double]"0.9999619126319885"
$ws.Cells.Item(42, 3).Value = [double]"2.537660066082026E-06"
$ws.Cells.Item(42, 4).Value = [double]"1"
$ws.Cells.Item(43, 1).Value = [double]"8.227544458350167E-05"
$ws.Cells.Item(43, 2).Value = [double]"0.9999809861183167"
$ws.Cells.Item(43, 3).Value = [double]"5.638419224851532E-07"
$ws.Cells.Item(43, 4).Value = [double]"1"
$ws.Cells.Item(44, 1).Value = [double]"2.270457298436668E-05"
$ws.Cells.Item(44, 2).Value = [double]"1"
$ws.Cells.Item(44, 3).Value = [double]"8.767211312488143E-08"
$ws.Cells.Item(44, 4).Value = [double]"1"
$ws.Cells.Item(45, 1).Value = [double]"0.0004820745089091361"
$ws.Cells.Item(45, 2).Value = [double]"0.9999428987503052"
$ws.Cells.Item(45, 3).Value = [double]"1.518641283837496E-07"
$ws.Cells.Item(45, 4).Value = [double]"1"
$ws.Cells.Item(46, 1).Value = [double]"0.0001530785812065005"
$ws.Cells.Item(46, 2).Value = [double]"0.9999428987503052"
$ws.Cells.Item(46, 3).Value = [double]"2.28959137871243E-07"
$ws.Cells.Item(46, 4).Value = [double]"1"
$ws.Cells.Item(47, 1).Value = [double]"0.0008504717843607068"
$ws.Cells.Item(47, 2).Value = [double]"0.9999238848686218"
$ws.Cells.Item(47, 3).Value = [double]"3.12010740799451E-07"
$ws.Cells.Item(47, 4).Value = [double]"1"
$ws.Cells.Item(48, 1).Value = [double]"0.000341077073244378"
$ws.Cells.Item(48, 2).Value = [double]"0.9999619126319885"
$ws.Cells.Item(48, 3).Value = [double]"4.823482981919369E-07"
$ws.Cells.Item(48, 4).Value = [double]"1"
$ws.Cells.Item(49, 1).Value = [double]"8.52015073178336E-05"
$ws.Cells.Item(49, 2).Value = [double]"0.9999428987503052"
$ws.Cells.Item(49, 3).Value = [double]"3.638583834231213E-08"
$ws.Cells.Item(49, 4).Value = [double]"1"
$ws.Cells.Item(50, 1).Value = [double]"7.497555634472519E-05"
$ws.Cells.Item(50, 2).Value = [double]"0.9999619126319885"
$ws.Cells.Item(50, 3).Value = [double]"2.746278227050425E-08"
$ws.Cells.Item(50, 4).Value = [double]"1"
$ws.Cells.Item(51, 1).Value = [double]"0.0002182240859838203"
$ws.Cells.Item(51, 2).Value = [double]"0.9999048113822937"
$ws.Cells.Item(51, 3).Value = [double]"0.0002710481639951468"
$ws.Cells.Item(51, 4).Value = [double]"1"
